# Apply inventory-consistency refresh ("finalizando inatividade de produtos")
# Updates estoque/data_estoque/estoque_depois/data_movimento (columns C, D, E, F)
# for the set of rows whose stock-check values changed, and backfills the
# F (data_movimento) cell for two rows that previously lacked one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row, new C/E (estoque / estoque_depois), new D (data_estoque), new F (data_movimento)
$data = @(
    @(33, 3153, 45818.80148693069, 45818.40175925926),
    @(38, 77, 45818.80148693405, 45818.46909722222),
    @(56, 121, 45818.8014869416, 45818.65388888889),
    @(82, 26, 45818.80148693439, 45818.46909722222),
    @(85, 0, 45818.80148693339, 45818.42027777778),
    @(87, 35, 45818.80148695138, 45818.7134375),
    @(91, 131, 45818.80148695299, 45818.71952546296),
    @(101, 580, 45818.80148694194, 45818.65388888889),
    @(106, 110, 45818.80148695203, 45818.71346064815),
    @(120, 16, 45818.80148693471, 45818.46909722222),
    @(124, 258, 45818.80148694228, 45818.65388888889),
    @(185, 136, 45818.80148693504, 45818.46909722222),
    @(195, 4, 45818.80148694695, 45818.66138888889),
    @(255, -1, 45818.8014869426, 45818.65388888889),
    @(274, 441, 45818.80148694293, 45818.65388888889),
    @(283, 218, 45818.80148693536, 45818.46909722222),
    @(291, 657, 45818.80148694325, 45818.65388888889),
    @(295, 288, 45818.80148693572, 45818.46909722222),
    @(309, 1115, 45818.80148693619, 45818.46909722222),
    @(326, 23, 45818.80148693661, 45818.46909722222),
    @(342, 212, 45818.80148694019, 45818.53353009259),
    @(385, 178, 45818.80148694391, 45818.65388888889),
    @(418, 148, 45818.80148693696, 45818.46909722222),
    @(510, 279, 45818.80148693114, 45818.40175925926),
    @(538, 624, 45818.8014869473, 45818.704375),
    @(542, 110, 45818.80148694778, 45818.704375),
    @(569, 8, 45818.80148694661, 45818.66094907407),
    @(570, 2622, 45818.80148693728, 45818.46909722222),
    @(631, 63, 45818.80148694831, 45818.704375),
    @(683, 32, 45818.80148694866, 45818.704375),
    @(691, 167, 45818.80148693761, 45818.46909722222),
    @(763, 112, 45818.80148694425, 45818.65388888889),
    @(852, -4, 45818.80148693794, 45818.46909722222),
    @(883, 70, 45818.80148694458, 45818.65388888889),
    @(891, -80, 45818.80148694057, 45818.53353009259),
    @(899, 0, 45818.80148694901, 45818.704375),
    @(923, 328, 45818.80148693373, 45818.42027777778),
    @(967, 2, 45818.80148695403, 45818.704375),
    @(972, 33, 45818.80148694938, 45818.704375),
    @(1002, 135, 45818.80148695243, 45818.71809027778),
    @(1062, 617, 45818.80148693828, 45818.46909722222),
    @(1122, 25, 45818.80148695171, 45818.7134375),
    @(1133, 82, 45818.80148694971, 45818.704375),
    @(1155, 13, 45818.80148695005, 45818.704375),
    @(1175, 61, 45818.80148695334, 45818.7203125),
    @(1187, 134, 45818.80148694491, 45818.65388888889),
    @(1222, 42, 45818.8014869315, 45818.40175925926),
    @(1228, 0, 45818.80148693191, 45818.40175925926),
    @(1247, 81, 45818.80148693859, 45818.46909722222),
    @(1253, -28, 45818.80148693892, 45818.46909722222),
    @(1283, -2, 45818.8014869504, 45818.704375),
    @(1325, -2, 45818.80148695073, 45818.704375),
    @(1342, 1040, 45818.80148693925, 45818.46909722222),
    @(1509, 139, 45818.80148695368, 45818.72087962963),
    @(1576, 26, 45818.80148693236, 45818.40175925926),
    @(1597, 5653, 45818.80148694525, 45818.65388888889),
    @(1600, 36, 45818.80148693273, 45818.40175925926),
    @(1932, 17, 45818.80148694559, 45818.65388888889),
    @(1934, 21, 45818.80148694627, 45818.65571759259),
    @(2022, 31, 45818.80148693304, 45818.40175925926),
    @(2023, 20, 45818.80148694593, 45818.65388888889),
    @(2081, -12, 45818.80148695105, 45818.704375),
    @(2416, 227, 45818.80148693962, 45818.46909722222)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $rownum = $row[0]
    $cval = $row[1]
    $dval = $row[2]
    $fval = $row[3]

    $ws.Cells.Item($rownum, 3).Value = $cval   # C: estoque
    $ws.Cells.Item($rownum, 4).Value = $dval   # D: data_estoque
    $ws.Cells.Item($rownum, 5).Value = $cval   # E: estoque_depois
    $ws.Cells.Item($rownum, 6).Value = $fval   # F: data_movimento
}

# Two rows (2561, 2562) gain a data_movimento (F) value where none existed
# before; give the new cell the same date/time number format used by the
# rest of column F/D so it matches the workbook's existing date style.
$dateFormat = $ws.Range("D2561").NumberFormat()

$ws.Range("F2561").NumberFormat = $dateFormat
$ws.Range("F2561").Value = 45818.63560185185

$ws.Range("F2562").NumberFormat = $dateFormat
$ws.Range("F2562").Value = 45818.63516203704
